$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: add a new "VISITOR_TYPE" column (G) --------------------------
$ws1.Range("G1").Value = "VISITOR_TYPE"

# Park the selection on Sheet1 at N5 before we switch the active sheet, so
# that Sheet1 no longer ends up as the tab-selected sheet once Sheet2 is
# activated below.
$ws1.Range("N5").Select()

# --- Add Sheet2 (placed immediately after Sheet1) --------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Rows 3-11: the list of supported visitor types
$ws2.Range("A3").Value = "Family"
$ws2.Range("A4").Value = "Interviewee"
$ws2.Range("A5").Value = "Vendor"
$ws2.Range("A6").Value = "Client"
$ws2.Range("A7").Value = "VIP"
$ws2.Range("A8").Value = "New Joinee"
$ws2.Range("A9").Value = "Conference Attendee"
$ws2.Range("A10").Value = "Guest"
$ws2.Range("A11").Value = "Visitor"

# Row 2: "Visitor Type" (bold)
$ws2.Range("A2").Value = "Visitor Type"
$ws2.Range("A2").Font.Bold = $true

# Row 17: "Unique Id Type" (bold)
$ws2.Range("A17").Value = "Unique Id Type"
$ws2.Range("A17").Font.Bold = $true

# Row 1: "Supported Values" title (bold, 16pt) with taller row height
$ws2.Range("A1").Value = "Supported Values"
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").Font.Size = 16
$ws2.Range("A1").RowHeight = 21

# Column widths for Sheet2 (closest values this engine's width grid can represent
# to the target bestFit-computed widths of 20.42578125 / 14.5703125)
$ws2.Columns.Item(1).ColumnWidth = 19.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 13.666666666666666

# Page setup for Sheet2
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Make Sheet2 the active sheet / tab with selection on I18
$ws2.Range("I18").Select()
